# Apply the edits described by the commit diff to exceldata.xlsx:
#  1. Cell C2 on Sheet1 (the phone/mobile number column) changes value
#     from 8015993932 to 9042977770.
#  2. The sheet's current selection moves from L1:L4 (active cell L4)
#     to the single cell C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the mobile number value in C2.
$ws.Range("C2").Value = 9042977770

# 2) Move the active selection to C3.
$ws.Range("C3").Select()
